$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" timestamp shown in the title cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 07:46"

# Full refreshed country table (columns A:H), rows 4-200.
# Country order changed (new entries interleaved / resorted) and the
# per-country statistics (Casos totales, Nuevos casos, Casos activos,
# Recuperados, Casos criticos, Muertes hoy, Muertes) were refreshed.
$data = @(
  ,@("China", 81218, 47, 73650, 4287, 1399, 4, 3281)
  ,@("Italia", 69176, 0, 8326, 54030, 3393, 0, 6820)
  ,@("Estados Unidos", 54916, 35, 379, 53753, 1175, 4, 784)
  ,@("España", 42058, 0, 3794, 35273, 2636, 0, 2991)
  ,@("Alemania", 32991, 0, 3290, 29542, 23, 0, 159)
  ,@("Iran", 24811, 0, 8913, 13964, 0, 0, 1934)
  ,@("Francia", 22304, 0, 3281, 17923, 2516, 0, 1100)
  ,@("Suiza", 9877, 0, 131, 9624, 141, 0, 122)
  ,@("Corea del Sur", 9137, 100, 3730, 5281, 59, 6, 126)
  ,@("Reino Unido", 8077, 0, 135, 7520, 20, 0, 422)
  ,@("Paises Bajos", 5560, 0, 2, 5282, 435, 0, 276)
  ,@("Austria", 5283, 0, 9, 5246, 22, 0, 28)
  ,@("Belgica", 4269, 0, 461, 3686, 381, 0, 122)
  ,@("Noruega", 2866, 0, 6, 2848, 44, 0, 12)
  ,@("Canada", 2792, 0, 112, 2654, 1, 0, 26)
  ,@("Portugal", 2362, 0, 22, 2307, 48, 0, 33)
  ,@("Australia", 2317, 0, 118, 2191, 11, 0, 8)
  ,@("Suecia", 2299, 0, 16, 2243, 143, 0, 40)
  ,@("Brasil", 2247, 0, 2, 2199, 18, 0, 46)
  ,@("Israel", 1930, 0, 53, 1874, 34, 0, 3)
  ,@("Turquia", 1872, 0, 0, 1828, 0, 0, 44)
  ,@("Malasia", 1796, 172, 183, 1597, 64, 0, 16)
  ,@("Dinamarca", 1591, 0, 1, 1558, 69, 0, 32)
  ,@("Chequia", 1394, 0, 10, 1381, 19, 0, 3)
  ,@("Irlanda", 1329, 0, 5, 1317, 29, 0, 7)
  ,@("Japon", 1193, 0, 285, 865, 54, 0, 43)
  ,@("Luxemburgo", 1099, 0, 6, 1085, 3, 0, 8)
  ,@("Ecuador", 1082, 0, 3, 1052, 2, 0, 27)
  ,@("Pakistan", 1000, 28, 19, 974, 5, 0, 7)
  ,@("Tailandia", 934, 107, 52, 878, 7, 0, 4)
  ,@("Chile", 922, 0, 17, 903, 7, 0, 2)
  ,@("Polonia", 901, 0, 1, 890, 3, 0, 10)
  ,@("Rumania", 794, 0, 79, 703, 15, 0, 12)
  ,@("Finlandia", 792, 0, 10, 781, 11, 0, 1)
  ,@("Arabia Saudita", 767, 0, 28, 738, 0, 0, 1)
  ,@("Grecia", 743, 0, 29, 694, 35, 0, 20)
  ,@("Crucero", 712, 0, 587, 115, 15, 0, 10)
  ,@("Sudafrica", 709, 155, 4, 705, 2, 0, 0)
  ,@("Indonesia", 686, 0, 30, 601, 0, 0, 55)
  ,@("Islandia", 648, 0, 51, 595, 13, 0, 2)
  ,@("India", 562, 26, 40, 512, 0, 0, 10)
  ,@("Singapur", 558, 0, 156, 400, 14, 0, 2)
  ,@("Filipinas", 552, 0, 20, 497, 1, 0, 35)
  ,@("Catar", 526, 0, 41, 485, 6, 0, 0)
  ,@("Rusia", 495, 0, 22, 472, 8, 0, 1)
  ,@("Eslovenia", 480, 0, 3, 473, 12, 0, 4)
  ,@("Panama", 443, 0, 1, 436, 33, 0, 6)
  ,@("Peru", 416, 0, 1, 408, 9, 0, 7)
  ,@("Egipto", 402, 0, 80, 302, 0, 0, 20)
  ,@("Barein", 392, 0, 177, 212, 2, 0, 3)
  ,@("Argentina", 387, 0, 52, 329, 0, 0, 6)
  ,@("Hong Kong", 387, 0, 102, 281, 4, 0, 4)
  ,@("Croacia", 382, 0, 5, 376, 6, 0, 1)
  ,@("Colombia", 378, 0, 6, 369, 0, 0, 3)
  ,@("Estonia", 369, 0, 7, 362, 4, 0, 0)
  ,@("Mexico", 367, 0, 4, 359, 1, 0, 4)
  ,@("Libano", 318, 0, 8, 306, 4, 0, 4)
  ,@("Irak", 316, 0, 75, 214, 0, 0, 27)
  ,@("Republica Dominicana", 312, 0, 3, 303, 0, 0, 6)
  ,@("Serbia", 303, 0, 15, 285, 21, 0, 3)
  ,@("Armenia", 265, 16, 16, 249, 6, 0, 0)
  ,@("Argelia", 264, 0, 65, 180, 0, 0, 19)
  ,@("Emiratos Arabes Unidos", 248, 0, 45, 201, 2, 0, 2)
  ,@("Hungria", 226, 39, 21, 195, 6, 1, 10)
  ,@("Bulgaria", 220, 2, 3, 214, 8, 0, 3)
  ,@("Taiwan", 216, 0, 29, 185, 0, 0, 2)
  ,@("Lituania", 209, 0, 1, 206, 1, 0, 2)
  ,@("Nueva Zelanda", 205, 0, 22, 183, 0, 0, 0)
  ,@("Eslovaquia", 204, 0, 7, 197, 2, 0, 0)
  ,@("Letonia", 197, 0, 1, 196, 0, 0, 0)
  ,@("Kuwait", 191, 0, 39, 152, 5, 0, 0)
  ,@("Uruguay", 189, 0, 0, 189, 3, 0, 0)
  ,@("San Marino", 187, 0, 4, 162, 12, 0, 21)
  ,@("Costa Rica", 177, 0, 2, 173, 4, 0, 2)
  ,@("Marruecos", 170, 0, 6, 159, 1, 0, 5)
  ,@("Bosnia y Herzegovina", 168, 0, 2, 163, 1, 0, 3)
  ,@("Principado de Andorra", 164, 0, 1, 162, 7, 0, 1)
  ,@("Jordania", 154, 0, 1, 153, 0, 0, 0)
  ,@("Republica de Macedonia", 148, 0, 1, 145, 1, 0, 2)
  ,@("Vietnam", 134, 0, 17, 117, 3, 0, 0)
  ,@("Moldavia", 125, 0, 2, 122, 10, 0, 1)
  ,@("Republica de Chipre", 124, 0, 3, 118, 3, 0, 3)
  ,@("Albania", 123, 0, 10, 108, 2, 0, 5)
  ,@("Islas Feroe", 122, 0, 33, 89, 0, 0, 0)
  ,@("Tunez", 114, 0, 1, 109, 11, 0, 4)
  ,@("Burkina Faso", 114, 0, 7, 103, 0, 0, 4)
  ,@("Malta", 110, 0, 2, 108, 1, 0, 0)
  ,@("Brunei", 104, 0, 2, 102, 2, 0, 0)
  ,@("Sri Lanka", 102, 0, 2, 100, 2, 0, 0)
  ,@("Ucrania", 102, 0, 1, 98, 0, 0, 3)
  ,@("Reunion", 94, 0, 1, 93, 0, 0, 0)
  ,@("Camboya", 93, 2, 4, 89, 1, 0, 0)
  ,@("Venezuela", 91, 7, 15, 76, 2, 0, 0)
  ,@("Azerbaiyan", 87, 0, 10, 76, 6, 0, 1)
  ,@("Senegal", 86, 0, 8, 78, 0, 0, 0)
  ,@("Oman", 84, 0, 17, 67, 0, 0, 0)
  ,@("Bielorrusia", 81, 0, 22, 59, 0, 0, 0)
  ,@("Kazajistan", 79, 7, 0, 79, 0, 0, 0)
  ,@("Afganistan", 74, 0, 1, 72, 0, 0, 1)
  ,@("Guadalupe", 73, 0, 0, 72, 4, 0, 1)
  ,@("Costa de Marfil", 73, 0, 2, 71, 0, 0, 0)
  ,@("Georgia", 70, 0, 9, 61, 1, 0, 0)
  ,@("Camerun", 66, 0, 2, 63, 0, 0, 1)
  ,@("Estado de Palestina", 60, 0, 16, 44, 0, 0, 0)
  ,@("Trinidad yTobago", 57, 0, 0, 57, 0, 0, 0)
  ,@("Martinica", 57, 0, 0, 56, 7, 0, 1)
  ,@("Uzbekistan", 55, 5, 0, 55, 4, 0, 0)
  ,@("Ghana", 53, 0, 0, 51, 0, 0, 2)
  ,@("Liechtenstein", 51, 0, 0, 51, 0, 0, 0)
  ,@("Cuba", 48, 0, 1, 46, 2, 0, 1)
  ,@("Montenegro", 47, 0, 0, 46, 0, 0, 1)
  ,@("Consejo Danes para los Refugiados", 45, 0, 0, 43, 0, 0, 2)
  ,@("Nigeria", 44, 0, 2, 41, 0, 0, 1)
  ,@("Kirguistan", 42, 0, 0, 42, 0, 0, 0)
  ,@("Mauricio", 42, 0, 0, 40, 1, 0, 2)
  ,@("Ruanda", 40, 0, 0, 40, 0, 0, 0)
  ,@("Puerto Rico", 39, 0, 1, 36, 0, 0, 2)
  ,@("Banglades", 39, 0, 5, 30, 0, 0, 4)
  ,@("Paraguay", 37, 10, 0, 34, 1, 1, 3)
  ,@("Mayotte", 36, 0, 0, 36, 0, 0, 0)
  ,@("Honduras", 36, 6, 0, 36, 0, 0, 0)
  ,@("Bolivia", 32, 3, 0, 32, 0, 0, 0)
  ,@("Guam", 32, 0, 0, 31, 0, 0, 1)
  ,@("Macao", 29, 0, 10, 19, 0, 0, 0)
  ,@("Kenia", 25, 0, 0, 25, 0, 0, 0)
  ,@("Polinesia Francesa", 25, 0, 0, 25, 0, 0, 0)
  ,@("Jamaica", 25, 4, 2, 22, 0, 0, 1)
  ,@("Isla de Man", 23, 0, 0, 23, 0, 0, 0)
  ,@("Monaco", 23, 0, 1, 22, 0, 0, 0)
  ,@("Guayana Francesa", 23, 0, 6, 17, 0, 0, 0)
  ,@("Guatemala", 21, 0, 0, 20, 0, 0, 1)
  ,@("Togo", 20, 0, 1, 19, 0, 0, 0)
  ,@("Barbados", 18, 0, 0, 18, 0, 0, 0)
  ,@("Islas Virgenes de los Estados Unidos", 17, 0, 0, 17, 0, 0, 0)
  ,@("Madagascar", 17, 0, 0, 17, 0, 0, 0)
  ,@("Aruba", 17, 0, 1, 16, 0, 0, 0)
  ,@("Gibraltar", 15, 0, 5, 10, 0, 0, 0)
  ,@("Maldivas", 13, 0, 5, 8, 0, 0, 0)
  ,@("Nueva Caledonia", 12, 2, 0, 12, 0, 0, 0)
  ,@("Etiopia", 12, 0, 0, 12, 0, 0, 0)
  ,@("Tanzania", 12, 0, 0, 12, 0, 0, 0)
  ,@("Mongolia", 10, 0, 0, 10, 0, 0, 0)
  ,@("El Salvador", 9, 4, 0, 9, 0, 0, 0)
  ,@("Guinea Ecuatorial", 9, 0, 0, 9, 0, 0, 0)
  ,@("Uganda", 9, 0, 0, 9, 0, 0, 0)
  ,@("San Martin (Parte Francesa)", 8, 0, 0, 8, 0, 0, 0)
  ,@("Haiti", 7, 0, 0, 7, 0, 0, 0)
  ,@("Seychelles", 7, 0, 0, 7, 0, 0, 0)
  ,@("Surinam", 7, 0, 0, 7, 0, 0, 0)
  ,@("Dominica", 7, 0, 0, 7, 0, 0, 0)
  ,@("Namibia", 7, 0, 2, 5, 0, 0, 0)
  ,@("Bermudas", 6, 0, 0, 6, 0, 0, 0)
  ,@("Benin", 6, 0, 0, 6, 0, 0, 0)
  ,@("Islas Caimanes", 6, 0, 0, 5, 0, 0, 1)
  ,@("Curazao", 6, 0, 0, 5, 0, 0, 1)
  ,@("Gabon", 6, 0, 0, 5, 0, 0, 1)
  ,@("Fiyi", 5, 1, 0, 5, 0, 0, 0)
  ,@("Guyana", 5, 0, 0, 4, 0, 0, 1)
  ,@("Bahamas", 5, 0, 1, 4, 0, 0, 0)
  ,@("Groenlandia", 5, 0, 2, 3, 0, 0, 0)
  ,@("Suazilandia", 4, 0, 0, 4, 0, 0, 0)
  ,@("Guinea", 4, 0, 0, 4, 0, 0, 0)
  ,@("Congo", 4, 0, 0, 4, 0, 0, 0)
  ,@("Santa Sede", 4, 0, 0, 4, 0, 0, 0)
  ,@("Liberia", 3, 0, 0, 3, 0, 0, 0)
  ,@("Mozambique", 3, 0, 0, 3, 0, 0, 0)
  ,@("San Bartolome", 3, 0, 0, 3, 0, 0, 0)
  ,@("Republica del Chad", 3, 0, 0, 3, 0, 0, 0)
  ,@("Republica de Yibuti", 3, 0, 0, 3, 0, 0, 0)
  ,@("Zambia", 3, 0, 0, 3, 0, 0, 0)
  ,@("Birmania", 3, 0, 0, 3, 0, 0, 0)
  ,@("Antigua y Barbuda", 3, 0, 0, 3, 0, 0, 0)
  ,@("Niger", 3, 0, 0, 3, 0, 0, 0)
  ,@("Santa Lucia", 3, 0, 0, 3, 0, 0, 0)
  ,@("Angola", 3, 0, 0, 3, 0, 0, 0)
  ,@("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0)
  ,@("Cabo Verde", 3, 0, 0, 2, 0, 0, 1)
  ,@("Zimbabue", 3, 0, 0, 2, 0, 0, 1)
  ,@("Sudan", 3, 0, 0, 2, 0, 0, 1)
  ,@("Gambia", 3, 0, 0, 2, 0, 0, 1)
  ,@("Butan", 2, 0, 0, 2, 0, 0, 0)
  ,@("Laos", 2, 0, 0, 2, 0, 0, 0)
  ,@("Nicaragua", 2, 0, 0, 2, 0, 0, 0)
  ,@("Mauritania", 2, 0, 0, 2, 0, 0, 0)
  ,@("San Martin (Parte Holandesa)", 2, 0, 0, 2, 0, 0, 0)
  ,@("Nepal", 2, 0, 1, 1, 0, 0, 0)
  ,@("Eritrea", 1, 0, 0, 1, 0, 0, 0)
  ,@("Libia", 1, 0, 0, 1, 0, 0, 0)
  ,@("Islas Turcas y Caicos", 1, 0, 0, 1, 0, 0, 0)
  ,@("Montserrat", 1, 0, 0, 1, 0, 0, 0)
  ,@("Belice", 1, 0, 0, 1, 0, 0, 0)
  ,@("San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0)
  ,@("Granada", 1, 0, 0, 1, 0, 0, 0)
  ,@("Siria", 1, 0, 0, 1, 0, 0, 0)
  ,@("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0)
  ,@("Somalia", 1, 0, 0, 1, 0, 0, 0)
  ,@("Timor Oriental", 1, 0, 0, 1, 0, 0, 0)
)

$r = 4
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
  $ws.Cells.Item($r, 8).Value = $row[7]
  $r++
}
